# Update the "想去人数" (number of interested people) values in the
# "展览" and "全部类型" worksheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 428
$ws1.Range("F5").Value  = 1730
$ws1.Range("F7").Value  = 2178
$ws1.Range("F11").Value = 4917
$ws1.Range("F12").Value = 10
$ws1.Range("F14").Value = 304
$ws1.Range("F17").Value = 182
$ws1.Range("F20").Value = 121
$ws1.Range("F21").Value = 3862
$ws1.Range("F22").Value = 711
$ws1.Range("F23").Value = 663
$ws1.Range("F26").Value = 104
$ws1.Range("F27").Value = 118
$ws1.Range("F28").Value = 23
$ws1.Range("F30").Value = 91
$ws1.Range("F34").Value = 946
$ws1.Range("F35").Value = 2457

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 428
$ws4.Range("F5").Value  = 1730
$ws4.Range("F7").Value  = 2178
$ws4.Range("F11").Value = 4917
$ws4.Range("F12").Value = 10
$ws4.Range("F14").Value = 304
$ws4.Range("F17").Value = 182
$ws4.Range("F20").Value = 121
$ws4.Range("F21").Value = 3862
$ws4.Range("F22").Value = 711
$ws4.Range("F23").Value = 663
$ws4.Range("F26").Value = 104
$ws4.Range("F27").Value = 118
$ws4.Range("F28").Value = 23
$ws4.Range("F30").Value = 91
$ws4.Range("F35").Value = 946
$ws4.Range("F36").Value = 2457
